$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$xlPasteFormats = -4122

# A new "2020" column (N) is being added to the table, to the right of the
# existing "2019" column (M). Carry over M's cell formatting (number format,
# font, borders, etc.) into N for the three data rows, then fill in the
# 2020 figures.
$ws.Range("M3").Copy() | Out-Null
$ws.Range("N3").PasteSpecial($xlPasteFormats) | Out-Null

$ws.Range("M4").Copy() | Out-Null
$ws.Range("N4").PasteSpecial($xlPasteFormats) | Out-Null

$ws.Range("M5").Copy() | Out-Null
$ws.Range("N5").PasteSpecial($xlPasteFormats) | Out-Null

$ws.Range("N3").Value = 2020
$ws.Range("N4").Value = 15
$ws.Range("N5").Value = 1308.3

# Leave the selection where the author left it after the edit.
$ws.Range("N6").Select() | Out-Null
